# relatorio_promar.xlsx — roll the monthly report from "Julho" to "Agosto":
# reset the summary totals to 0 and drop this month's (not-yet-collected)
# detail rows for each "Top N" breakdown, keeping only each section's
# header + column-label row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the month and zero out the headline totals.
$ws.Range("B2").Value = "Agosto"
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0

# Remove the per-section detail rows (bottom-up so unprocessed row numbers
# don't shift while we work).
$ws.Rows("30:33").Delete()   # Idades Mais Visitadas - detail rows
$ws.Rows("23:26").Delete()   # Notas Mais Dadas - detail rows
$ws.Rows("18:19").Delete()   # Cidades com Melhor Desempenho - detail rows
$ws.Rows("9:14").Delete()    # Respostas Mais Acertadas - detail rows
